$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Expand the table "Tabelle1" by three rows (A1:H26 -> A1:H29).
$null = $tbl.ListRows.Add()
$null = $tbl.ListRows.Add()
$null = $tbl.ListRows.Add()

# Copy the formatting of the last pre-existing data row (row 26) down onto the
# three brand-new rows so they pick up the same borders / fonts as the rest
# of the table before we layer the wrap-text formatting on top.
$ws.Range("A26:H26").Copy()
$ws.Range("A27:H29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new rows in the same cell-by-cell order the values were
# originally authored in, so newly-introduced shared strings line up with
# the target workbook.
# Row 27: "Sel_Crit_compare" summary row
$ws.Range("C27").Value = "Sel_Crit_compare"
$ws.Range("D27").Value = "Zusammenfassung der expliziten Werte für das Selektieren von RNA-Dependence"
$ws.Range("E27").Value = "delta_mu, Ctrl/Rnase_heightprop, overlap, Ctrl/Rnase_heightprop_compare"
$ws.Range("F27").Value = "5990 x 6"
$ws.Range("G27").Value = "Chunk 29"
$ws.Range("H27").Value = "Sel_Crit"

# Chunk numbers for rows 28-29
$ws.Range("G28").Value = "Chunk 31"
$ws.Range("G29").Value = "Chunk 31"

# Selection-criteria labels
$ws.Range("C28").Value = "RDPs"
$ws.Range("C29").Value = "non-RDPs"

# Dimensions
$ws.Range("F28").Value = "706 x 6"
$ws.Range("F29").Value = "5284 x 6"

# Descriptions (entered last)
$ws.Range("D28").Value = "Alle Proteine, bei denen alle drei RDP-Selektionskriterien erfüllt sind"
$ws.Range("D29").Value = "Alle Proteine, bei denen mind. Ein RDP-Selektionskriterium nicht erfüllt ist"

# Variablen & Herkunft columns reuse already-existing text
$ws.Range("E28").Value = "delta_mu, Ctrl/Rnase_heightprop, overlap, Ctrl/Rnase_heightprop_compare"
$ws.Range("E29").Value = "delta_mu, Ctrl/Rnase_heightprop, overlap, Ctrl/Rnase_heightprop_compare"
$ws.Range("H28").Value = "Sel_Crit_compare"
$ws.Range("H29").Value = "Sel_Crit_compare"

# "#" column numbering
$ws.Range("A27").Value = 13
$ws.Range("A28").Value = 14
$ws.Range("A29").Value = 14

# The three new rows get wrap-text applied to every column except "#" and
# "Alter Name" (matching the long-text presentation already used for row
# 26's Beschreibung/Variablen, which also leaves A/B unwrapped).
$ws.Range("C27:H29").WrapText = $true

# New row heights shrink to fit the (shorter) wrapped text in these rows.
$ws.Range("27:29").RowHeight = 29

# Keep the table's hidden AutoFilter helper name in sync with the new extent.
$filterName = $ws.Names.Item("_FilterDatabase")
$filterName.RefersTo = "=Tabelle1!`$C`$1:`$C`$29"

$wb.Save()
